$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# --- Paragraph: "Specific divisions within NCU-F..." ---
# "...have dozens of area owners, which creates routing challenges." ->
# "...have dozens of area owners, creating routing challenges."
$find1 = "have dozens of area owners, which creates routing challenges."
$replace1 = "have dozens of area owners, creating routing challenges."
Replace-Text $find1 $replace1

# "...Incident Management software, such as  PagerDuty, or a similar..." ->
# "...Incident Management software like PagerDuty, or a similar..."
$find2 = "Incident Management software, such as  PagerDuty, or a similar"
$replace2 = "Incident Management software like PagerDuty, or a similar"
Replace-Text $find2 $replace2

# --- Paragraph: "Suppose the engineering team cannot mitigate the issue..." ---
# "...Those third-party providers experience similar economic constraints, ... service team." ->
# "...Those third-party providers have similar economic constraints, ... service team.  However,
#  ... Lastly, reporting must inform the executive leadership of any business impact to approve
#  future investments into the problem."
$find3 = "Those third-party providers experience similar economic constraints, which forces them into these stacked pyramid structures.  This consistency includes knowledge databases, support channels, and customer access to the service team."
$replace3 = "Those third-party providers have similar economic constraints, which forces them into these stacked pyramid structures.  This consistency includes knowledge databases, support channels, and customer access to the service team.  However, there can be specific scenarios that are not resolvable.  In these situations, both the engineering and support program managers need to agree on an appropriate response.  Responses can include adding to the backlog, proposing workarounds, among other stopgaps.  Lastly, reporting must inform the executive leadership of any business impact to approve future investments into the problem."
Replace-Text $find3 $replace3

Write-Output "Edits applied successfully."
